# Shopizer Signin Page - Test Case sheet cleanup
# The author removed a block of 7 blank placeholder rows (rows 52-58) that
# separated the SC_01_TC_10 block from SC_01_TC_13/14/15, and renumbered the
# following three test cases (SC_01_TC_13 -> SC_01_TC_11, SC_01_TC_14 ->
# SC_01_TC_12, SC_01_TC_15 -> SC_01_TC_13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Case")

# Delete the 7 blank rows (52-58); everything below shifts up by 7.
$ws.Rows("52:58").Delete()

# Rename the Tc_id values for the three test cases that moved up.
$ws.Range("A53").Value = "SC_01_TC_11"
$ws.Range("A57").Value = "SC_01_TC_12"
$ws.Range("A61").Value = "SC_01_TC_13"

# Update the sheet view to match the saved state in the workbook.
$ws.Activate()
$ws.Range("A61").Select()
